# Individual Scenario KDMA scores (UK DATA) (#396)
# * update headers
# * var defs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns for the Moral Judgement (MJ) narrative / non-narrative
#     KDMA text scores, right before the existing MJ_KDMA_Sim column (Q). ---
$ws.Columns("Q:R").Insert()

# --- Insert two new columns for the Ingroup Bias (IO) narrative / non-narrative
#     KDMA text scores, right before the existing IO_KDMA_Sim column (now U). ---
$ws.Columns("U:V").Insert()

# ----------------------------------------------------------------------------
# Row 1 - Variable names (header row)
# ----------------------------------------------------------------------------
$ws.Range("Q1").Value = "MJ_KDMA_Text_Narr"
$ws.Range("R1").Value = "MJ_KDMA_Text_NonNarr"
$ws.Range("U1").Value = "IO_KDMA_Text_Narr"
$ws.Range("V1").Value = "IO_KDMA_Text_NonNarr"

# ----------------------------------------------------------------------------
# Row 2 - Variable Group
# ----------------------------------------------------------------------------
$ws.Range("Q2").Value = "Attribute Assessment"
$ws.Range("R2").Value = "Attribute Assessment"
$ws.Range("U2").Value = "Attribute Assessment"
$ws.Range("V2").Value = "Attribute Assessment"

# ----------------------------------------------------------------------------
# Row 3 - Description
# ----------------------------------------------------------------------------
$ws.Range("Q3").Value = "Moral Judgement KDMA measurement from narrative Adept text scenario only (MJ5)"
$ws.Range("R3").Value = "Moral Judgement KDMA measurement from non narrative Adept text scenario only (MJ1)"
$ws.Range("U3").Value = "Ingroup Bias KDMA measurement from narrative Adept text scenario only (MJ5)"
$ws.Range("V3").Value = "Ingroup Bias KDMA measurement from non narrative Adept text scenario only (IO1)"

# ----------------------------------------------------------------------------
# Row 4 - Labels
# ----------------------------------------------------------------------------
$ws.Range("Q4").Value = "Number"
$ws.Range("R4").Value = "Number"
$ws.Range("U4").Value = "Number"
$ws.Range("V4").Value = "Number"

# ----------------------------------------------------------------------------
# Row 5 - Calculation
# ----------------------------------------------------------------------------
$ws.Range("Q5").Value = "-"
$ws.Range("R5").Value = "-"
$ws.Range("U5").Value = "-"
$ws.Range("V5").Value = "-"

# D5:F5 no longer hold an empty-string shared value - clear them to true blanks
# and drop their quote-prefixed format (copy the plain format used by the rest
# of row 5, e.g. G5, then blank the values).
$ws.Range("G5").Copy()
$ws.Range("D5:F5").PasteSpecial(-4122)
$ws.Range("D5:F5").Value = ""

# ----------------------------------------------------------------------------
# Row 6 - Source
# ----------------------------------------------------------------------------
$ws.Range("Q6").Value = "From TA1 Server"
$ws.Range("R6").Value = "From TA1 Server"
$ws.Range("U6").Value = "From TA1 Server"
$ws.Range("V6").Value = "From TA1 Server"

# --- Re-fit every column to its (new) widest content, matching the rest of
#     the sheet's auto-sized columns. ---
$ws.Columns("A:W").AutoFit()
